# Fixed error in computation of a simulated moment
# Updates the "simulation" column (D) values on the "data" worksheet to
# reflect the corrected simulated moments.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("D5").Value  = 2.485535850014337
$ws.Range("D6").Value  = 0.07175346706947595
$ws.Range("D7").Value  = -0.2569518751681636
$ws.Range("D8").Value  = 0.2463665638042128
$ws.Range("D9").Value  = 2.482289296325269
$ws.Range("D10").Value = 0.2846642908342636
$ws.Range("D11").Value = 2.454420695467284
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 0.6261929753734357
$ws.Range("D22").Value = 0.4671769704044489
$ws.Range("D23").Value = 0.1929157180451532
